# Update the generalized_cf.xlsx performance results with the latest
# algorithm run's mean/std values (re-upload of refreshed results).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.006641094635598557
$ws.Range("C2").Value = 0.007491279236003888
$ws.Range("D2").Value = 0.007066186935801223
$ws.Range("E2").Value = 0.0004250923002026654

$ws.Range("B3").Value = 0.1911037891268534
$ws.Range("C3").Value = 0.216887417218543
$ws.Range("D3").Value = 0.2039956031726982
$ws.Range("E3").Value = 0.01289181404584483

$ws.Range("B4").Value = 0.01283611818081222
$ws.Range("C4").Value = 0.01448233928472721
$ws.Range("D4").Value = 0.01365922873276972
$ws.Range("E4").Value = 0.0008231105519574981
